# edit.ps1
# Applies the case14.xlsx edit:
#  - reorders sheets: bus, branch, generator, generatorcost (new), basePower, info
#  - adds a new "generatorcost" worksheet with cost-model data
#  - sets the active tab to the new "generatorcost" sheet

$wb = $excel.ActiveWorkbook

# --- 1. Reorder existing sheets: swap "generator" and "branch" --------------
$generator = $wb.Worksheets.Item("generator")
$branch    = $wb.Worksheets.Item("branch")
$generator.Move($null, $branch)          # generator moves to just after branch

# --- 2. Insert the new "generatorcost" sheet before "basePower" ------------
$basePower = $wb.Worksheets.Item("basePower")
$costSheet = $wb.Worksheets.Add($basePower)
$costSheet.Name = "generatorcost"

# --- 3. Populate the new sheet ----------------------------------------------
# Row 1: headers
$costSheet.Cells.Item(1,1).Value = "Cost Model"
$costSheet.Cells.Item(1,2).Value = "Cost"
$costSheet.Cells.Item(1,3).Value = "Cost"
$costSheet.Cells.Item(1,4).Value = "Cost Model"
$costSheet.Cells.Item(1,5).Value = "Cost Model"
$costSheet.Cells.Item(1,6).Value = "Cost Model"
$costSheet.Cells.Item(1,7).Value = "Cost Model"

# Row 2: sub-headers
$costSheet.Cells.Item(2,1).Value = "Piecewise(1), Polynomial(2)"
$costSheet.Cells.Item(2,2).Value = "Startup [currency]"
$costSheet.Cells.Item(2,3).Value = "Shutdown [currency]"
$costSheet.Cells.Item(2,4).Value = "Number of Data Points"
$costSheet.Cells.Item(2,5).Value = "Coefficient c2"
$costSheet.Cells.Item(2,6).Value = "Coefficient c1"
$costSheet.Cells.Item(2,7).Value = "Coefficient c0"

# Rows 3-7: data
$data = @(
    @(2, 0, 0, 3, 0.043029259899999998, 20, 0),
    @(2, 0, 0, 3, 0.25,                  20, 0),
    @(2, 0, 0, 3, 0.01,                  40, 0),
    @(2, 0, 0, 3, 0.01,                  40, 0),
    @(2, 0, 0, 3, 0.01,                  40, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 3
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $costSheet.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}

# --- 4. Restore the generator sheet's remembered selection (C9) -------------
# (re-fetch by name: the worksheet handle can go stale across a Move())
$generator = $wb.Worksheets.Item("generator")
$generator.Activate()
$generator.Range("C9").Select()

# --- 5. Make the new sheet the active tab -----------------------------------
$costSheet = $wb.Worksheets.Item("generatorcost")
$costSheet.Activate()
$costSheet.Range("I9").Select()
